# Swap the data (columns B:AB) between pairs of adjacent rows.
# Column A (the running index) stays put; every other field - id, teams,
# odds, results, etc. - moves from one row to the other, effectively
# swapping the two match records while keeping their position (A) fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(17, 18),
    @(41, 42),
    @(58, 59),
    @(108, 109),
    @(135, 136),
    @(151, 152),
    @(161, 162)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AB$r1")
    $rng2 = $ws.Range("B$r2`:AB$r2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}
